$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export rolled forward by one day: the oldest day (2025-11-07)
# is dropped, and every following day's row shifts up to take its place.
$ws.Rows.Item(2).Delete()
